$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.764.11"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "2.438.85"
$ws.Range("E3").Value = "  -1.74%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.89"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.18"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.507"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.169"
$ws.Range("E9").Value = "  +6.99%  "
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.331"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("E12").Value = "  -5.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000179"
$ws.Range("E13").Value = "  +3.95%  "
$ws.Range("D14").Value = "68.736.25"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "2.890.28"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.34"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").Value = "2.445.07"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.59"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.80"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.03"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.94"
$ws.Range("E22").Value = "  +1.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.79"
$ws.Range("E24").Value = "  -1.71%  "
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").Value = "2.569.44"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.30"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").Value = "0.0₃0824"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.17"
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "432.10"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  -2.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "160.13"
$ws.Range("E35").Value = "  +0.60%  "
$ws.Range("E36").Value = "  -0.17%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.00"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.107"
$ws.Range("E39").Value = "  -1.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.300"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.51"
$ws.Range("E41").Value = "  +1.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.37"
$ws.Range("E42").Value = "  -2.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.07"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("E44").Value = "  -0.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.33"
$ws.Range("E45").Value = "  -1.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.51"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("E48").Value = "  -1.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.558"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("E51").Value = "  +0.33%  "
